# Insert two new weekly records (rows 489-490) into the "Macroferia Regional
# de Talca - Plátano" data block, pushing the existing rows 489-577 down to
# 491-579 (dimension grows from A1:T577 to A1:T579).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 489 (Excel shifts rows 489.. down).
$ws.Rows.Item(489).Insert()
$ws.Rows.Item(489).Insert()

# --- New row 489: "Pintón" record -----------------------------------------
$ws.Range("A489").Value = 5
$ws.Range("B489").Value = "Macroferia Regional de Talca"
$ws.Range("C489").Value = "Maule"
$ws.Range("D489").Value = 44637
$ws.Range("E489").Value = 7
$ws.Range("F489").Value = "Fruta"
$ws.Range("G489").Value = 100108
$ws.Range("H489").Value = "Tropicales y subtropicales"
$ws.Range("I489").Value = 100108006
$ws.Range("J489").Value = "Plátano"
$ws.Range("K489").Value = "Sin especificar"
$ws.Range("L489").Value = "Pintón"
$ws.Range("M489").Value = 500
$ws.Range("N489").Value = 17000
$ws.Range("O489").Value = 17000
$ws.Range("P489").Value = 17000
$ws.Range("Q489").Value = "$/caja 20 kilos"
$ws.Range("R489").Value = "Ecuador"
$ws.Range("S489").Value = 850
$ws.Range("T489").Value = 20

# --- New row 490: "Primera Pintón" record ----------------------------------
$ws.Range("A490").Value = 5
$ws.Range("B490").Value = "Macroferia Regional de Talca"
$ws.Range("C490").Value = "Maule"
$ws.Range("D490").Value = 44637
$ws.Range("E490").Value = 7
$ws.Range("F490").Value = "Fruta"
$ws.Range("G490").Value = 100108
$ws.Range("H490").Value = "Tropicales y subtropicales"
$ws.Range("I490").Value = 100108006
$ws.Range("J490").Value = "Plátano"
$ws.Range("K490").Value = "Sin especificar"
$ws.Range("L490").Value = "Primera Pintón"
$ws.Range("M490").Value = 300
$ws.Range("N490").Value = 18000
$ws.Range("O490").Value = 18000
$ws.Range("P490").Value = 18000
$ws.Range("Q490").Value = "$/caja 20 kilos"
$ws.Range("R490").Value = "Ecuador"
$ws.Range("S490").Value = 900
$ws.Range("T490").Value = 20
